# Apply the record-pair swaps described by the diff.
# Each pair of rows below had their species/location data swapped with
# each other (the row numbers stay put, but the content that belongs to
# each observation moves to the other row).  Only the specific cells that
# actually differ are touched, so every other cell (dates, booleans,
# lookups, etc.) is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 40 / 42 -----------------------------------------------------
$ws.Range("A40").Value = 130803045
$ws.Range("M40").Value = "färska spår"
$ws.Range("Q40").Value = 424768
$ws.Range("R40").Value = 6712134

$ws.Range("A42").Value = 130803059
$ws.Range("M42").Value = "äldre spår"
$ws.Range("Q42").Value = 424858
$ws.Range("R42").Value = 6712137

# --- Rows 47 / 48 -------------------------------------------------------
$ws.Range("A47").Value = 130803052
$ws.Range("Q47").Value = 424773
$ws.Range("R47").Value = 6712133

$ws.Range("A48").Value = 130803049
$ws.Range("Q48").Value = 424771
$ws.Range("R48").Value = 6712443

# --- Rows 56 / 57 -------------------------------------------------------
$ws.Range("A56").Value = 130803048
$ws.Range("Q56").Value = 424801
$ws.Range("R56").Value = 6712325

$ws.Range("A57").Value = 130803057
$ws.Range("Q57").Value = 424935
$ws.Range("R57").Value = 6712079

# --- Rows 61 / 62 (Tretåig hackspett <-> Mörk kolflarnlav) -------------
$ws.Range("A61").Value = 130848907
$ws.Range("B61").Value = 57884
$ws.Range("E61").Value = 100109
$ws.Range("F61").Value = "Tretåig hackspett"
$ws.Range("G61").Value = "Picoides tridactylus"
$ws.Range("H61").Value = "(Linnaeus, 1758)"
$ws.Range("K61").Value = ""
$ws.Range("L61").Value = ""
$ws.Range("M61").Value = "äldre spår"
$ws.Range("N61").Value = ""
$ws.Range("Q61").Value = 424588
$ws.Range("R61").Value = 6712316
$ws.Range("AC61").Value = "Ringhack"

$ws.Range("A62").Value = 130848917
$ws.Range("B62").Value = 79002
$ws.Range("E62").Value = 228912
$ws.Range("F62").Value = "Mörk kolflarnlav"
$ws.Range("G62").Value = "Carbonicola myrmecina"
$ws.Range("H62").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("K62").Value = ""
$ws.Range("L62").Value = ""
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = ""
$ws.Range("Q62").Value = 424590
$ws.Range("R62").Value = 6712294
$ws.Range("AC62").Value = ""

# --- Rows 65 / 66 -------------------------------------------------------
$ws.Range("A65").Value = 130848929
$ws.Range("B65").Value = 79001
$ws.Range("E65").Value = 6446
$ws.Range("F65").Value = "Kolflarnlav"
$ws.Range("G65").Value = "Carbonicola anthracophila"
$ws.Range("H65").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q65").Value = 424592
$ws.Range("R65").Value = 6712413

$ws.Range("A66").Value = 130848922
$ws.Range("B66").Value = 81229
$ws.Range("E66").Value = 1049
$ws.Range("F66").Value = "Kortskaftad ärgspik"
$ws.Range("G66").Value = "Microcalicium ahlneri"
$ws.Range("H66").Value = "Tibell"
$ws.Range("Q66").Value = 424576
$ws.Range("R66").Value = 6712311
